# "Update matrix and list"
#
# The document contains a single 3-column relationship matrix table.
# This script:
#   1. Updates the row ("Cal Or Info" | "Problematic TTDs" | "1:1") so its
#      second cell reads "Calibration Racks" instead of "Problematic TTDs".
#   2. Removes the now-redundant row
#      ("Cal Or Info" | "Calibration Racks" | "1:N") entirely.
#   3. Updates the row ("SO Info" | "TTD Rack Info" | "1:1") so its second
#      cell reads "Tube Test Devices" instead of "TTD Rack Info".

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. "Problematic TTDs" -> "Calibration Racks" (row: Cal Or Info | ... | 1:1)
$cell = $t.Rows.Item(26).Cells.Item(2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("Problematic TTDs", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Calibration Racks", 1)

# --- 2. Delete the row "Cal Or Info | Calibration Racks | 1:N"
$t.Rows.Item(27).Delete()

# --- 3. "TTD Rack Info" -> "Tube Test Devices" (row: SO Info | ... | 1:1)
$cell2 = $t.Rows.Item(27).Cells.Item(2)
$r2 = $d.Range($cell2.Range.Start, $cell2.Range.End)
$r2.Find.Execute("TTD Rack Info", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "Tube Test Devices", 1)
